# Branches-id.xlsx update: append newly-onboarded branches to the list
# (id/name table on Sheet1) and tidy up the sheet the way Excel would
# after you type new rows in and resize the columns to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the three new branch rows under the existing 36 rows of data ---
$newRows = @(
    @(182, "فرع بريدة 1"),
    @(183, "فرع حائل 2"),
    @(196, "فرع راحتي خميس مشيط 2")
)

$r = 37
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# --- Resize columns A and B to fit the (now longer) content ---
$ws.Columns.Item(1).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(2).ColumnWidth = 21.5

# --- Turn the AutoFilter off now that the list has been edited ---
$ws.AutoFilterMode = $false
$wb.Names.Item("_xlnm._FilterDatabase").Delete()

# --- Reset the page margins back to the Excel defaults ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# --- Leave the selection where the editor last clicked ---
[void]$ws.Range("E21").Select()
